$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths
# Note: the ColumnWidth property snaps to a 1/6-character grid internally,
# so the inputs below are chosen as the closest values that round to the
# intended stored widths (3.140625, 2.140625 and 5.7109375 characters).
$ws.Columns.Item(3).ColumnWidth = 2.3333333333333335
$ws.Columns.Item(5).ColumnWidth = 2.3333333333333335
$ws.Columns.Item(6).ColumnWidth = 1.3333333333333333
$ws.Columns.Item(7).ColumnWidth = 1.3333333333333333
$ws.Columns.Item(11).ColumnWidth = 4.833333333333333
$ws.Columns.Item(12).ColumnWidth = 4.833333333333333

# Update cell values in row 1
$ws.Range("B1").Value = 4
$ws.Range("C1").Value = 23
$ws.Range("D1").Value = 12
$ws.Range("E1").Value = 17
$ws.Range("F1").Value = 8
$ws.Range("G1").Value = 2
$ws.Range("H1").Value = 7
$ws.Range("I1").Value = 22
$ws.Range("J1").Value = 30
$ws.Range("K1").Value = 0.020999999999999998
$ws.Range("L1").Value = 0.007000000000000001
$ws.Range("M1").Value = 0.079000000000000001
$ws.Range("N1").Value = 0.076999999999999999
